$d = $word.ActiveDocument

# Move to the very end of the document (inside the last, currently-empty
# paragraph that already precedes the sectPr) and push a fresh paragraph
# mark after it so the existing empty <w:p/> is left untouched.
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()

# The new last paragraph is where we inject the OOXML fragment for the
# guide text. Collapse to its start so InsertXML adds content *before*
# this (still empty) paragraph instead of overwriting it.
$newLast = $d.Paragraphs.Last
$insertionPoint = $newLast.Range
$insertionPoint.Collapse(1)  # wdCollapseStart

$fragment = @'
<w:p/>
<w:p><w:r><w:t>Guide:</w:t></w:r></w:p>
<w:p/>
<w:p><w:r><w:t xml:space="preserve">The timed command variants are there to stop the sprite before water hazards as we can’t just lift our fingers off a key to stop it.  </w:t></w:r><w:r><w:t>Here’s a stab at a guide:</w:t></w:r></w:p>
<w:p/>
<w:p><w:r><w:t>Level 1:</w:t></w:r></w:p>
<w:p><w:r><w:t xml:space="preserve">Command sequence: ‘right’ + ‘jump’ – you need to say ‘jump’ directly after right so that he can pick up the coin. </w:t></w:r></w:p>
<w:p/>
<w:p><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">Level 2: </w:t></w:r></w:p>
<w:p><w:r><w:t xml:space="preserve">When you encounter ladders here is the sequence: </w:t></w:r></w:p>
<w:p/>
<w:p><w:r><w:t>‘up’ + ‘&lt;</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>wiggle|jiggle|turbo|blast|boggle</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>&gt;’  - all but ‘boggle’ will jump to the right, boggle is not needed as much but is all you need ladders where the platform at the top is to the left.</w:t></w:r></w:p>
<w:p/>
<w:p><w:r><w:t xml:space="preserve">Level 3: </w:t></w:r></w:p>
<w:p><w:r><w:t xml:space="preserve">First water hazard. </w:t></w:r></w:p>
<w:p><w:r><w:t>Command sequence: ‘right’ + ‘jiggle’   - should bring you to the edge of the water hazard at a dead stop – then use turbo to get over (you might miss the coin though)</w:t></w:r></w:p>
<w:p/>
<w:p/>
<w:p><w:r><w:t xml:space="preserve">Level 4: </w:t></w:r></w:p>
<w:p><w:r><w:t xml:space="preserve">Two water hazards, this is the trickiest one. Progress carefully.  For jumping over stuff you can use ‘jiggle’ or </w:t></w:r><w:r><w:t xml:space="preserve">even ‘turbo’ instead of wiggle as the former are timed and it may be safer. </w:t></w:r></w:p>
'@

$insertionPoint.InsertXML($fragment)
